# Auto-generated edit script: updates the cryptos price/volume table
# to the latest scrape (GitHub Actions refresh), cell-by-cell, to match
# the authoritative diff. Each text-like numeric string (e.g. "1.00",
# "0.0930") is entered with a leading apostrophe so Excel stores it as
# literal text instead of coercing it to a Number (which would drop
# trailing zeros / introduce floating point noise), then the style is
# reset back to Normal so no stray number-format/quote-prefix survives.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.036.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.82%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.304.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'110.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.45%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'313.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.85%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.82%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.35%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.27%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'44.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.16%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.21%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'8.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.45%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +17.55%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -0.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.07%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'2.642.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.38%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.350.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.95%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'43.034.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.31%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.26%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'76.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.39%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.81%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +7.24%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'257.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +11.24%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.72%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'11.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.81%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.11%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'39.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.58%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.46%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'22.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +5.26%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'173.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.85%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0903"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.30%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.01%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.87%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.51%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -7.57%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0377"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.83%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.91%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'71.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.03%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +10.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.51%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.09%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'12.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.50%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'5.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'108.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.56%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +4.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.37%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'ordi"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'70.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.33%  "
$ws.Range("E51").Style = "Normal"
